$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.888.83'
$ws.Range('E2').Value = '  -2.50%  '
$ws.Range('D3').Value = '3.422.11'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''582.36'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.11%  '
$ws.Range('D6').Value = '''173.69'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.63%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '''0.591'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -4.11%  '
$ws.Range('D9').Value = '3.421.82'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('D10').Value = '''0.130'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -8.40%  '
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('E12').Value = '  -5.30%  '
$ws.Range('D13').Value = '4.017.37'
$ws.Range('E13').Value = '  -1.39%  '
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').Value = '''29.90'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -7.06%  '
$ws.Range('D16').Value = '65.928.35'
$ws.Range('E16').Value = '  -2.41%  '
$ws.Range('E17').Value = '  -4.27%  '
$ws.Range('D18').Value = '3.423.98'
$ws.Range('E18').Value = '  -1.41%  '
$ws.Range('D19').Value = '''5.86'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.59%  '
$ws.Range('D20').Value = '''13.69'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.42%  '
$ws.Range('D21').Value = '''365.66'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -7.47%  '
$ws.Range('E22').Value = '  -3.49%  '
$ws.Range('D23').Value = '''0.998'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  -1.80%  '
$ws.Range('D25').Value = '''71.44'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.70%  '
$ws.Range('D26').Value = '''0.524'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.00%  '
$ws.Range('D27').Value = '''0.0000117'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.92%  '
$ws.Range('D28').Value = '''9.65'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -7.09%  '
$ws.Range('E29').Value = '  +0.54%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').Value = '''23.87'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.23%  '
$ws.Range('D32').Value = '''5.75'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.38%  '
$ws.Range('D33').Value = '''1.98'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.01%  '
$ws.Range('D34').Value = '''1.00'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -8.69%  '
$ws.Range('D36').Value = '''7.01'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.57%  '
$ws.Range('D37').Value = '''1.54'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.86%  '
$ws.Range('D38').Value = '''159.62'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.83%  '
$ws.Range('D39').Value = '''28.88'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +9.84%  '
$ws.Range('E40').Value = '  -1.98%  '
$ws.Range('D41').Value = '''1.78'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.75%  '
$ws.Range('D42').Value = '''2.55'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -11.10%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.700.57'
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '''4.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -6.34%  '
$ws.Range('D45').Value = '''6.30'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.23%  '
$ws.Range('D46').Value = '''0.0677'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -6.00%  '
$ws.Range('D47').Value = '''39.90'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').Value = '''24.02'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -8.98%  '
$ws.Range('D49').Value = '''0.0288'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.77%  '
$ws.Range('D50').Value = '''304.07'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -6.98%  '
$ws.Range('D51').Value = '''0.811'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.21%  '
